$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore C10 (row 10, column C) from 18 to 1, per the commit's "Restored
# from revision" of the rules table (min value for rule R20).
$ws.Range("C10").Value = 1
